# Correção das notas do fórum para matc65 em 2021.2
# For every student row whose "nota_view" (column J) is non-zero, zero out
# the daily view flags (columns B-H), the total_views (column I) and the
# nota_view (column J) itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $jVal = $ws.Cells.Item($r, 10).Value2
    if ($jVal -ne 0 -and $jVal -ne $null) {
        $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 10)).Value = 0
    }
}
